$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update row 3 and row 4 in column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1213
$wsExhibition.Range("F4").Value = 2686

# Sheet "全部类型" (All types): update row 5 and row 6 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1213
$wsAll.Range("F6").Value = 2686
